$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '57.769.50'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -5.60%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.907.09'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -3.50%  '

$ws.Range("E4").Value = '  +0.10%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '550.33'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.12%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '122.90'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -4.61%  '

$ws.Range("E7").Value = '  +0.13%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '2.903.25'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -3.55%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.495'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.41%  '

$ws.Range("E10").Value = '  -7.00%  '

$ws.Range("E11").Value = '  -7.98%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.439'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +2.14%  '

$ws.Range("E13").Value = '  -4.59%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '32.34'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.48%  '

$ws.Range("E15").Value = '  +1.37%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.384.45'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -3.54%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.906.49'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -3.64%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '6.59'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +6.09%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '57.760.67'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -5.61%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '409.41'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -6.56%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '12.91'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.78%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.672'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.68%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.85'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -3.99%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '12.88'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +2.81%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '77.20'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.29%  '

$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("E27").Value = '  +0.00%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.46'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.70%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.25'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.97%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.94'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +3.31%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.06'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.52%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '24.72'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.31%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0983'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +4.74%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.912'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -4.25%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '5.38'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.07%  '

$ws.Range("E36").Value = '  -11.48%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '48.14'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -3.91%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '8.48'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +9.77%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0₃0621'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -8.22%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0345'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -5.06%  '

$ws.Range("E41").Value = '  -1.47%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.627.41'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.72%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '360.57'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.58%  '

$ws.Range("E44").Value = '  -1.12%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.02%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '120.18'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.68%  '

$ws.Range("E47").Value = '  -2.78%  '

$ws.Range("E48").Value = '  +0.94%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.94'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.45%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '22.81'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -3.34%  '

$ws.Range("E51").Value = '  -2.80%  '
